$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.015.82'
$ws.Range('E2').Value = '  -3.83%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.959.30'
$ws.Range('E3').Value = '  -5.73%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.005'
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '327.12'
$ws.Range('E5').Value = '  -3.20%  '
$ws.Range('E6').Value = '  +0.08%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4956'
$ws.Range('E7').Value = '  -5.24%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.4195'
$ws.Range('E8').Value = '  -3.63%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '52.72'
$ws.Range('E9').Value = '  -3.91%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.09229'
$ws.Range('E10').Value = '  -0.97%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.095'
$ws.Range('E11').Value = '  -6.20%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '22.81'
$ws.Range('E12').Value = '  -6.59%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.984.55'
$ws.Range('E13').Value = '  -0.45%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.435'
$ws.Range('E14').Value = '  -5.87%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.826'
$ws.Range('E15').Value = '  -7.14%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.007'
$ws.Range('E16').Value = '  +0.25%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '91.26'
$ws.Range('E17').Value = '  -8.93%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.00001097'
$ws.Range('E18').Value = '  -5.05%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06690'
$ws.Range('E19').Value = '  -0.06%  '
$ws.Range('E20').Value = '  -7.62%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.004'
$ws.Range('E21').Value = '  +0.27%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.928'
$ws.Range('E22').Value = '  -5.53%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '29.067.91'
$ws.Range('E23').Value = '  -3.74%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '11.98'
$ws.Range('E24').Value = '  -3.04%  '
$ws.Range('E25').Value = '  -2.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.190.96'
$ws.Range('E26').Value = '  -3.91%  '
$ws.Range('E27').Value = '  -5.26%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '155.62'
$ws.Range('E28').Value = '  -3.97%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '6.276'
$ws.Range('E29').Value = '  -7.41%  '
$ws.Range('E30').Value = '  -9.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '126.08'
$ws.Range('E31').Value = '  -5.24%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.040'
$ws.Range('E32').Value = '  -7.55%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.09812'
$ws.Range('E33').Value = '  -5.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.519'
$ws.Range('E34').Value = '  -8.01%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.812'
$ws.Range('E35').Value = '  -6.73%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '3.673'
$ws.Range('E36').Value = '  -5.78%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02420'
$ws.Range('E37').Value = '  -6.63%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.318'
$ws.Range('E38').Value = '  -0.36%  '
$ws.Range('B39').Value = 'FraxShare'
$ws.Range('C39').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '8.971'
$ws.Range('E39').Value = '  -8.35%  '
$ws.Range('B40').Value = 'Hedera'
$ws.Range('C40').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.06341'
$ws.Range('E40').Value = '  -5.36%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.6432'
$ws.Range('E41').Value = '  -6.94%  '
$ws.Range('E42').Value = '  -8.70%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1962'
$ws.Range('E43').Value = '  -10.56%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.365'
$ws.Range('E45').Value = '  +4.56%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.6196'
$ws.Range('E46').Value = '  -7.43%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '13.30'
$ws.Range('E47').Value = '  -6.13%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '2.188'
$ws.Range('E48').Value = '  -6.67%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.464'
$ws.Range('E49').Value = '  -4.17%  '
$ws.Range('E50').Value = '  -5.72%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.07028'
$ws.Range('E51').Value = '  -2.49%  '
